$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.784.14"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "3.463.54"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'585.83"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("D6").Value = "'178.50"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.97%  "
$ws.Range("D7").Value = "'0.631"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +6.54%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "3.465.48"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("E10").Value = "  +1.69%  "
$ws.Range("D11").Value = "'6.95"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").Value = "4.067.24"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").Value = "'30.19"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "66.623.11"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("D18").Value = "3.513.38"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "'13.91"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("D21").Value = "'372.98"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.57%  "
$ws.Range("D22").Value = "'7.68"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("D23").Value = "'73.60"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.21%  "
$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D24").Value = "'0.0000128"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +7.87%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("D27").Value = "'9.98"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("E28").Value = "  +2.80%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("D31").Value = "'2.00"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("D32").Value = "'23.71"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.40%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("E35").Value = "  -2.72%  "
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("D37").Value = "'162.81"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.06%  "
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("D39").Value = "'28.04"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.88%  "
$ws.Range("E40").Value = "  +2.87%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'2.60"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.86%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'4.51"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.44%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.763.23"
$ws.Range("E43").Value = "  +3.57%  "
$ws.Range("D44").Value = "'6.49"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.63%  "
$ws.Range("D45").Value = "'0.0698"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("D46").Value = "'25.45"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.37%  "
$ws.Range("D47").Value = "'341.54"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +9.34%  "
$ws.Range("D48").Value = "'40.08"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("D49").Value = "'0.0289"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").Value = "'0.106"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.72%  "
$ws.Range("D51").Value = "'31.82"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.77%  "
